# Update column F ("dSF") values on the active sheet to reflect the
# repulled/recalculated data, per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -1
    5  = 0
    6  = -1
    15 = 4
    17 = 1
    20 = -3
    25 = 2
    26 = -1
    32 = -3
    38 = 1
    41 = 3
    48 = -1
    49 = 2
    53 = -1
    56 = 2
    57 = -2
    65 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
